# Fruta / hortaliza, semanal
# Insert a new weekly price-record row before the current row 433,
# shifting all subsequent records (old rows 433-504) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 433:504 down to 434:505, creating a blank row 433.
$ws.Rows("433:433").Insert()

# Populate the newly inserted row 433 with the new weekly record.
$ws.Range("A433").Value = 5
$ws.Range("B433").Value = "Macroferia Regional de Talca"
$ws.Range("C433").Value = "Maule"
$ws.Range("D433").Value = 44476
$ws.Range("E433").Value = 7
$ws.Range("F433").Value = "Fruta"
$ws.Range("G433").Value = 100104
$ws.Range("H433").Value = "Frutos de pepita"
$ws.Range("I433").Value = 100104002
$ws.Range("J433").Value = "Manzana"
$ws.Range("K433").Value = "Granny Smith"
$ws.Range("L433").Value = "Primera"
$ws.Range("M433").Value = 250
$ws.Range("N433").Value = 8000
$ws.Range("O433").Value = 8000
$ws.Range("P433").Value = 8000
$ws.Range("Q433").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R433").Value = "Región de O'Higgins"
$ws.Range("S433").Value = 533
$ws.Range("T433").Value = 15
